$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "24.693.16"
$cell.ClearFormats()
$ws.Range("E2").Value = "  +2.27%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.696.17"
$cell.ClearFormats()
$ws.Range("E3").Value = "  +1.35%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "316.81"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +2.31%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.ClearFormats()
$ws.Range("E6").Value = "  +0.38%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.3961"
$cell.ClearFormats()
$ws.Range("E7").Value = "  +1.90%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.4165"
$cell.ClearFormats()
$ws.Range("E8").Value = "  +3.90%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "1.509"
$cell.ClearFormats()
$ws.Range("E9").Value = "  +2.93%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.ClearFormats()
$ws.Range("E10").Value = "  +0.32%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "51.52"
$cell.ClearFormats()
$ws.Range("E11").Value = "  -4.78%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.08825"
$cell.ClearFormats()
$ws.Range("E12").Value = "  +1.75%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "7.216"
$cell.ClearFormats()
$ws.Range("E13").Value = "  +4.69%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "23.49"
$cell.ClearFormats()
$ws.Range("E14").Value = "  +3.15%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "8.177"
$cell.ClearFormats()
$ws.Range("E15").Value = "  +11.35%  "
$ws.Range("E16").Value = "  +0.87%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "1.698.81"
$cell.ClearFormats()
$ws.Range("E17").Value = "  +2.19%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "100.22"
$cell.ClearFormats()
$ws.Range("E18").Value = "  +0.88%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.07024"
$cell.ClearFormats()
$ws.Range("E19").Value = "  +1.25%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "19.78"
$cell.ClearFormats()
$ws.Range("E20").Value = "  +3.08%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "7.114"
$cell.ClearFormats()
$ws.Range("E21").Value = "  +8.18%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.9999"
$cell.ClearFormats()
$ws.Range("E22").Value = "  +0.43%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "14.32"
$cell.ClearFormats()
$ws.Range("E23").Value = "  +2.24%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "24.679.88"
$cell.ClearFormats()
$ws.Range("E24").Value = "  +2.19%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "3.119"
$cell.ClearFormats()
$ws.Range("E25").Value = "  +2.30%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.339"
$cell.ClearFormats()
$ws.Range("E26").Value = "  +1.41%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "22.85"
$cell.ClearFormats()
$ws.Range("E27").Value = "  +4.80%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "162.51"
$cell.ClearFormats()
$ws.Range("E28").Value = "  +1.10%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "137.32"
$cell.ClearFormats()
$ws.Range("E29").Value = "  +5.12%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "5.175"
$cell.ClearFormats()
$ws.Range("E30").Value = "  +1.09%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "7.460"
$cell.ClearFormats()
$ws.Range("E31").Value = "  +2.01%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "1.884.10"
$cell.ClearFormats()
$ws.Range("E32").Value = "  +1.90%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "1.085"
$cell.ClearFormats()
$ws.Range("E33").Value = "  -1.55%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.08631"
$cell.ClearFormats()
$ws.Range("E34").Value = "  -0.28%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "7.149"
$cell.ClearFormats()
$ws.Range("E35").Value = "  -0.36%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "11.55"
$cell.ClearFormats()
$ws.Range("E36").Value = "  +1.80%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.2757"
$cell.ClearFormats()
$ws.Range("E37").Value = "  +3.69%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "1.935"
$cell.ClearFormats()
$ws.Range("E38").Value = "  +0.19%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "14.52"
$cell.ClearFormats()
$ws.Range("E39").Value = "  -0.32%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.09248"
$cell.ClearFormats()
$ws.Range("E40").Value = "  +4.57%  "
$ws.Range("E41").Value = "  +7.66%  "
$ws.Range("E42").Value = "  +2.38%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.7682"
$cell.ClearFormats()
$ws.Range("E43").Value = "  +1.47%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "2.657"
$cell.ClearFormats()
$ws.Range("E44").Value = "  +9.60%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "15.94"
$cell.ClearFormats()
$ws.Range("E45").Value = "  +5.06%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.7204"
$cell.ClearFormats()
$ws.Range("E46").Value = "  +1.24%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "4.220"
$cell.ClearFormats()
$ws.Range("E47").Value = "  +2.75%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.9996"
$cell.ClearFormats()
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("E49").Value = "  +5.98%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "141.17"
$cell.ClearFormats()
$ws.Range("E50").Value = "  +1.26%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.07995"
$cell.ClearFormats()
$ws.Range("E51").Value = "  +1.61%  "
